{"js": "// Gate/radio-button form template cleanup:\n// Replace the raw merge placeholders \"${roleInspektur}\" and \"${inspektur}\"\n// (signature block) with Indonesian instructional placeholder text, since\n// these fields are no longer auto-filled by the generator.\n\n// 1) \"${roleInspektur}\" -> \"(Masukkan jabatan penandatangan)\"\n//    In the source document this placeholder is split across two runs\n//    (\"${roleInspektur\" + \"}\"), but Body.search() matches across run\n//    boundaries, and Range.insertText(..., \"Replace\") collapses the\n//    matched range into a single run carrying the first run's formatting\n//    (matching the target: one <w:r> with the original rFonts/sz/szCs).\nconst roleResults = context.document.body.search(\"${roleInspektur}\", {\n  matchCase: true,\n  matchWholeWord: false,\n  matchWildcards: false\n});\nroleResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < roleResults.items.length; i++) {\n  roleResults.items[i].insertText(\"(Masukkan jabatan penandatangan)\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"${inspektur}\" -> \"(Masukkan nama penandatangan)\"\nconst inspekturResults = context.document.body.search(\"${inspektur}\", {\n  matchCase: true,\n  matchWholeWord: false,\n  matchWildcards: false\n});\ninspekturResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < inspekturResults.items.length; i++) {\n  inspekturResults.items[i].insertText(\"(Masukkan nama penandatangan)\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Gate/radio-button form template cleanup:\n# Replace the raw merge placeholders \"${roleInspektur}\" and \"${inspektur}\"\n# (signature block) with Indonesian instructional placeholder text, since\n# these fields are no longer auto-filled by the generator.\n\n$d = $word.ActiveDocument\n\n# 1) \"${roleInspektur}\" -> \"(Masukkan jabatan penandatangan)\"\n#    In the source document this placeholder is split across two runs\n#    (\"${roleInspektur\" + \"}\"); Find/Replace matches across run\n#    boundaries and collapses the hit into a single run that keeps the\n#    first run's formatting (matching the target: one <w:r> with the\n#    original rFonts/sz/szCs).\n$rngRole = $d.Content\n$null = $rngRole.Find.Execute(\n  \"`${roleInspektur}\",\n  $false, $false, $false, $false, $false,\n  $true, 1, $false,\n  \"(Masukkan jabatan penandatangan)\",\n  2\n)\n\n# 2) \"${inspektur}\" -> \"(Masukkan nama penandatangan)\"\n$rngInspektur = $d.Content\n$null = $rngInspektur.Find.Execute(\n  \"`${inspektur}\",\n  $false, $false, $false, $false, $false,\n  $true, 1, $false,\n  \"(Masukkan nama penandatangan)\",\n  2\n)\n"}
